$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.924.37'
$ws.Range("E2").Value = '  +0.63%  '
$ws.Range("D3").Value = '2.920.70'
$ws.Range("E3").Value = '  +0.90%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '593.62'
$ws.Range("E5").Value = '  +1.51%  '
$ws.Range("D6").Value = '145.63'
$ws.Range("E6").Value = '  -0.33%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +0.96%  '
$ws.Range("D9").Value = '6.82'
$ws.Range("E9").Value = '  +2.20%  '
$ws.Range("E10").Value = '  +0.23%  '
$ws.Range("D11").Value = '0.439'
$ws.Range("E11").Value = '  -1.89%  '
$ws.Range("E12").Value = '  +0.65%  '
$ws.Range("D13").Value = '33.70'
$ws.Range("E13").Value = '  -0.53%  '
$ws.Range("E14").Value = '  -0.08%  '
$ws.Range("D15").Value = '3.401.82'
$ws.Range("E15").Value = '  +0.87%  '
$ws.Range("D16").Value = '60.936.18'
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("D17").Value = '6.71'
$ws.Range("E17").Value = '  -1.11%  '
$ws.Range("D18").Value = '2.912.36'
$ws.Range("E18").Value = '  +0.64%  '
$ws.Range("D19").Value = '431.01'
$ws.Range("E19").Value = '  +1.60%  '
$ws.Range("D20").Value = '13.36'
$ws.Range("E20").Value = '  -1.77%  '
$ws.Range("D21").Value = '0.683'
$ws.Range("E21").Value = '  +1.87%  '
$ws.Range("D22").Value = '7.07'
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Value = '81.62'
$ws.Range("E23").Value = '  +1.90%  '
$ws.Range("D24").Value = '11.00'
$ws.Range("E24").Value = '  +0.16%  '
$ws.Range("E25").Value = '  +0.72%  '
$ws.Range("D26").Value = '11.95'
$ws.Range("E26").Value = '  +1.01%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("E28").Value = '  +4.87%  '
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("D30").Value = '2.61'
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("D31").Value = '7.05'
$ws.Range("E31").Value = '  -2.89%  '
$ws.Range("D32").Value = '26.41'
$ws.Range("E32").Value = '  +0.15%  '
$ws.Range("E33").Value = '  +1.01%  '
$ws.Range("D34").Value = '0.0₃0851'
$ws.Range("E34").Value = '  +1.85%  '
$ws.Range("E35").Value = '  +1.06%  '
$ws.Range("D36").Value = '5.63'
$ws.Range("E36").Value = '  -0.43%  '
$ws.Range("D37").Value = '3.02'
$ws.Range("E37").Value = '  +2.86%  '
$ws.Range("E38").Value = '  +0.32%  '
$ws.Range("D39").Value = '1.98'
$ws.Range("D40").Value = '8.56'
$ws.Range("E40").Value = '  -1.59%  '
$ws.Range("E41").Value = '  -1.78%  '
$ws.Range("D42").Value = '40.38'
$ws.Range("E42").Value = '  -3.03%  '
$ws.Range("D43").Value = '373.55'
$ws.Range("E43").Value = '  +0.37%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("D45").Value = '2.704.16'
$ws.Range("E45").Value = '  +2.04%  '
$ws.Range("D46").Value = '130.88'
$ws.Range("E46").Value = '  -1.40%  '
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("D48").Value = '23.98'
$ws.Range("E48").Value = '  -5.00%  '
$ws.Range("E49").Value = '  +0.24%  '
$ws.Range("E50").Value = '  -3.29%  '
$ws.Range("E51").Value = '  +2.58%  '
